# Update the "想去人数" (F column) figures in the "展览" and "全部类型" sheets
# to match the newly published numbers.

$wb = $excel.ActiveWorkbook

# Common updates that are identical on both sheets (row -> new value)
$commonUpdates = @{
    2  = 826
    4  = 1149
    5  = 49
    6  = 12369
    7  = 53
    8  = 104
    9  = 494
    10 = 444
    11 = 1131
    12 = 910
    13 = 13617
    14 = 13819
    19 = 1034
    20 = 104
    23 = 4919
    24 = 221
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $commonUpdates.Keys) {
        $ws.Range("F$row").Value = $commonUpdates[$row]
    }

    # F22 differs slightly between the two sheets before the edit,
    # but both converge to 190 afterwards.
    $ws.Range("F22").Value = 190
}
